$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3: extend the "red cells" note
$ws.Range("A3").Value = "红色格子：填入你的实验数据，如本身自带数据请更改"

# A32: update the GitHub repo URL
$ws.Range("A32").Value = "Posted on https://github.com/Axolyz/fuck-nku-physics-experiments."

# A34: clear the closing quote/text entirely
$ws.Range("A34").ClearContents()

# Leave the final selection on A3, matching the saved view state
$ws.Range("A3").Select()

$wb.Save()
